# #5: cash & deposit done
# Fix the "存款" (deposits) sheet: give it proper column headers and
# append the standard trailing metadata columns (property_category,
# category, date, legislator_name, legislator_id, source_file, index)
# that the other sheets (land / car / insurance) already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---- fix the header row (row 1) ------------------------------------
# Row 1 used to be a stray duplicate of the first data row; turn it
# into real column headers, matching the other sheets' naming.
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"

# ---- add the new trailing header cells (G1:M1) ----------------------
# Copy the bold/bordered header style from an existing header cell so
# the new cells match the sheet's look.
$ws.Range("B1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)

$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# ---- row 11's total was stored as text "1000000"; make it numeric --
$ws.Range("F11").Value = 1000000

# The "date" column (I) holds a literal text value that looks like a
# date ("2012-02-10"); force text formatting up front so Excel's
# auto-conversion doesn't turn it into a date serial number.
$ws.Range("I2:I17").NumberFormat = "@"

# ---- populate the new G:M columns for each data row (2-17) ---------
# Values are the same constant metadata used on every other sheet of
# this workbook (property_category=deposit, category=normal, the
# filing date, legislator name/id, source file), plus the per-row
# "index" which mirrors column A's row id.
for ($r = 2; $r -le 17; $r++) {
    # copy the plain data-row style onto the new cells first
    $ws.Range("B" + $r).Copy()
    $ws.Range("G" + $r + ":M" + $r).PasteSpecial(-4122)
    # re-assert text formatting on the date cell (PasteSpecial above
    # would otherwise overwrite it with the copied General format)
    $ws.Cells.Item($r, 9).NumberFormat = "@"

    $ws.Cells.Item($r, 7).Value = "deposit"
    $ws.Cells.Item($r, 8).Value = "normal"
    $ws.Cells.Item($r, 9).Value = "2012-02-10"
    $ws.Cells.Item($r, 10).Value = "盧秀燕"
    $ws.Cells.Item($r, 11).Value = 869
    $ws.Cells.Item($r, 12).Value = "tmp61a71"
    $ws.Cells.Item($r, 13).Value = $ws.Cells.Item($r, 1).Value()
}
